{"js": "// The Somass watershed risk table has one column layout:\n// 0 Watershed | 1 LF (description) | 2 Rank | 3 Total Risk | 4 Current Risk | 5 Future Risk\n//\n// The commit \"Adjusted risk calc formula\" recalculated the risk scores for a\n// block of \"Rank 15\" (LF7, LF8, LF12, LF33, LF65) and \"Rank 20\" (LF19, LF25,\n// LF30, LF70) limiting factors, which bumped some Total Risk / Current Risk /\n// Future Risk values and re-ordered the four \"Rank 20\" rows (by their LF\n// description text) while nudging one of them to Rank 23.\n//\n// Apply the edit as a fixed set of per-cell text replacements addressed by\n// (row index, column index) in the single table, which is both exactly what\n// changed in the underlying OOXML and robust against re-serialization.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row, column, old value (for a sanity check), new value\nconst cellEdits = [\n  [15, 3, \"4\", \"6\"],\n  [15, 5, \"L\", \"M\"],\n  [16, 3, \"4\", \"6\"],\n  [16, 5, \"L\", \"M\"],\n  [17, 3, \"4\", \"6\"],\n  [17, 5, \"L\", \"M\"],\n  [18, 3, \"4\", \"6\"],\n  [18, 5, \"L\", \"M\"],\n  [19, 3, \"4\", \"6\"],\n  [19, 5, \"L\", \"M\"],\n  [20, 1, \"LF19: Mortality or fitness reduction due to early alevin emergence\", \"LF25: Mortality or fitness reduction due to lower quality spawning gravel\"],\n  [20, 3, \"2\", \"4\"],\n  [20, 4, \"VL\", \"L\"],\n  [21, 1, \"LF25: Mortality or fitness reduction due to lower quality spawning gravel\", \"LF30: Mortality or fitness reduction as a result of elevated predation\"],\n  [21, 3, \"2\", \"4\"],\n  [21, 5, \"VL\", \"L\"],\n  [22, 1, \"LF30: Mortality or fitness reduction as a result of elevated predation\", \"LF70: Mortality or fitness reduction due to negative effects of small population size - including inbreeding depression and gene flow\"],\n  [22, 3, \"2\", \"4\"],\n  [22, 5, \"VL\", \"L\"],\n  [23, 1, \"LF70: Mortality or fitness reduction due to negative effects of small population size - including inbreeding depression and gene flow\", \"LF19: Mortality or fitness reduction due to early alevin emergence\"],\n  [23, 2, \"20\", \"23\"],\n  [23, 4, \"L\", \"VL\"],\n  [23, 5, \"VL\", \"L\"],\n];\n\nconst cells = cellEdits.map(([row, col]) => table.getCell(row, col));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < cellEdits.length; i++) {\n  const [, , oldValue, newValue] = cellEdits[i];\n  const cell = cells[i];\n  // Only overwrite when the cell still holds the value we expect; this keeps\n  // the script a no-op (per cell) if it is ever re-run against an\n  // already-edited document instead of corrupting unrelated content.\n  if (cell.value === oldValue) {\n    cell.value = newValue;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The Somass watershed risk table has one column layout:\n# 1 Watershed | 2 LF (description) | 3 Rank | 4 Total Risk | 5 Current Risk | 6 Future Risk\n#\n# The commit \"Adjusted risk calc formula\" recalculated the risk scores for a\n# block of \"Rank 15\" (LF7, LF8, LF12, LF33, LF65) and \"Rank 20\" (LF19, LF25,\n# LF30, LF70) limiting factors, which bumped some Total Risk / Current Risk /\n# Future Risk values and re-ordered the four \"Rank 20\" rows (by their LF\n# description text) while nudging one of them to Rank 23.\n#\n# Apply the edit as a fixed set of per-cell text replacements addressed by\n# (row, column) -- 1-based, as Word's Table.Cell() expects -- in the single\n# table, which is both exactly what changed in the underlying OOXML and\n# robust against re-serialization.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$edits = @(\n  @{Row=16; Col=4; Old=\"4\"; New=\"6\"},\n  @{Row=16; Col=6; Old=\"L\"; New=\"M\"},\n  @{Row=17; Col=4; Old=\"4\"; New=\"6\"},\n  @{Row=17; Col=6; Old=\"L\"; New=\"M\"},\n  @{Row=18; Col=4; Old=\"4\"; New=\"6\"},\n  @{Row=18; Col=6; Old=\"L\"; New=\"M\"},\n  @{Row=19; Col=4; Old=\"4\"; New=\"6\"},\n  @{Row=19; Col=6; Old=\"L\"; New=\"M\"},\n  @{Row=20; Col=4; Old=\"4\"; New=\"6\"},\n  @{Row=20; Col=6; Old=\"L\"; New=\"M\"},\n  @{Row=21; Col=2; Old=\"LF19: Mortality or fitness reduction due to early alevin emergence\"; New=\"LF25: Mortality or fitness reduction due to lower quality spawning gravel\"},\n  @{Row=21; Col=4; Old=\"2\"; New=\"4\"},\n  @{Row=21; Col=5; Old=\"VL\"; New=\"L\"},\n  @{Row=22; Col=2; Old=\"LF25: Mortality or fitness reduction due to lower quality spawning gravel\"; New=\"LF30: Mortality or fitness reduction as a result of elevated predation\"},\n  @{Row=22; Col=4; Old=\"2\"; New=\"4\"},\n  @{Row=22; Col=6; Old=\"VL\"; New=\"L\"},\n  @{Row=23; Col=2; Old=\"LF30: Mortality or fitness reduction as a result of elevated predation\"; New=\"LF70: Mortality or fitness reduction due to negative effects of small population size - including inbreeding depression and gene flow\"},\n  @{Row=23; Col=4; Old=\"2\"; New=\"4\"},\n  @{Row=23; Col=6; Old=\"VL\"; New=\"L\"},\n  @{Row=24; Col=2; Old=\"LF70: Mortality or fitness reduction due to negative effects of small population size - including inbreeding depression and gene flow\"; New=\"LF19: Mortality or fitness reduction due to early alevin emergence\"},\n  @{Row=24; Col=3; Old=\"20\"; New=\"23\"},\n  @{Row=24; Col=5; Old=\"L\"; New=\"VL\"},\n  @{Row=24; Col=6; Old=\"VL\"; New=\"L\"}\n)\n\nforeach ($e in $edits) {\n  $cell = $t.Cell($e.Row, $e.Col)\n  # Cell.Range.Text carries a trailing cell-mark (CR + BEL); strip it before\n  # comparing so the guard matches the visible cell content.\n  $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n  if ($current -eq $e.Old) {\n    $cell.Range.Text = $e.New\n  }\n}\n"}
